$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells stay text (avoid Excel auto-converting numeric-looking
# strings like "211.46" or "0.513" into real numbers) by forcing a text number format
# before assigning the values, matching the original inline-string cell content.

# Row 19 and Row 20 swap ranking order: BitcoinCash <-> Dai, with updated price/volume data
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.22%  "

# Price/volume updates for remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.664.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.63"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.578.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.642.35"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.04"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.04"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0514"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.287.74"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.04%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +16.39%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.782"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.17"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.732.90"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.09"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.37"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.15%  "
